# adding task description to xlsx doc
#
# Summary of edits (see commit message "adding task description to xlsx doc"):
#   - "protocol" sheet: insert a new column C ("full_task_name" header,
#     "empathy intervention" values for the 4 task rows), shifting the
#     former C..H columns to D..I.
#   - "participants" sheet: apply the existing "Arial" cell style (the same
#     style already used on A1) across the B1:D6 data block, clear the
#     (accidental) style on A2, and drop the stale conditional-formatting
#     rule that used to flag blank cells in columns B:D.
#   - Selection/active-cell bookkeeping on both sheets is updated to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "protocol": insert the full_task_name column
# ---------------------------------------------------------------
$protocol = $wb.Worksheets.Item("protocol")

# Insert a new blank column before column C; everything from the old C
# onward shifts one column to the right (old C->D, D->E, ... G->H, H->I).
$protocol.Columns("C").Insert()

# New column inherits column B's width so it renders the same as the
# other un-"bestFit" columns around it.
$protocol.Columns("C").ColumnWidth = $protocol.Columns("B").ColumnWidth

# Header + the four per-run values for the newly created column.
$protocol.Range("C1").Value = "full_task_name"
$protocol.Range("C5:C8").Value = "empathy intervention"

# ---------------------------------------------------------------
# Sheet "participants": re-style the data columns + cleanup
# ---------------------------------------------------------------
$participants = $wb.Worksheets.Item("participants")

# Drop the old "flag blank cells in B:D" conditional formatting rule.
$participants.Range("B1:D1048576").FormatConditions.Delete()

# Give columns B:D (header + all data rows) the same look as column A's
# header cell (A1) by copying its format across - this reuses the
# existing style rather than fabricating a new one.
$participants.Range("A1").Copy()
$participants.Range("B1:D6").PasteSpecial(-4122)

# A2 had accidentally picked up that same style; put it back to the
# workbook default, matching its sibling cells A3:A6.
$participants.Range("A2").ClearFormats()

# Move the stored selection/active cell on "participants" ...
$participants.Range("E17").Select()

# ... then re-activate "protocol" (the originally/still active tab) and
# restore its selection, so tabSelected stays on "protocol".
$protocol.Activate()
$protocol.Range("C6:C8").Select()
